$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new value is numeric-looking text that must remain TEXT
# (matches source formatting, e.g. trailing zeros like "1.000").
# Temporarily mark as Text format so Excel does not reinterpret the
# string as a number, then restore the original (Normal/General) style
# so the saved file keeps the default style index.
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D51").NumberFormat = "@"

# Apply the updated values
$ws.Range("D2").Value = '26.983.58'
$ws.Range("E2").Value = '  -0.09%  '
$ws.Range("D3").Value = '1.874.28'
$ws.Range("E3").Value = '  +0.63%  '
$ws.Range("E4").Value = '  +0.17%  '
$ws.Range("D5").Value = '305.55'
$ws.Range("E5").Value = '  -0.11%  '
$ws.Range("E6").Value = '  +0.17%  '
$ws.Range("D7").Value = '0.5062'
$ws.Range("E7").Value = '  -0.05%  '
$ws.Range("D8").Value = '0.3662'
$ws.Range("E8").Value = '  -2.06%  '
$ws.Range("D9").Value = '0.07203'
$ws.Range("E9").Value = '  +0.90%  '
$ws.Range("D10").Value = '0.8947'
$ws.Range("E10").Value = '  +1.39%  '
$ws.Range("D11").Value = '20.72'
$ws.Range("E11").Value = '  +0.46%  '
$ws.Range("B12").Value = 'TRON'
$ws.Range("C12").Value = 'https://coinranking.com/coin/qUhEFk1I61atv+tron-trx'
$ws.Range("D12").Value = '0.07529'
$ws.Range("E12").Value = '  -0.49%  '
$ws.Range("B13").Value = 'WrappedEther'
$ws.Range("C13").Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range("D13").Value = '1.862.40'
$ws.Range("E13").Value = '  -0.25%  '
$ws.Range("D14").Value = '94.70'
$ws.Range("E14").Value = '  +6.11%  '
$ws.Range("E15").Value = '  -1.24%  '
$ws.Range("D16").Value = '1.000'
$ws.Range("E16").Value = '  +0.16%  '
$ws.Range("D17").Value = '0.000008548'
$ws.Range("E17").Value = '  +1.72%  '
$ws.Range("D18").Value = '14.25'
$ws.Range("E18").Value = '  +1.49%  '
$ws.Range("D19").Value = '1.000'
$ws.Range("E19").Value = '  +0.23%  '
$ws.Range("D20").Value = '27.024.10'
$ws.Range("E20").Value = '  -0.10%  '
$ws.Range("E21").Value = '  -0.01%  '
$ws.Range("D22").Value = '2.089.92'
$ws.Range("E22").Value = '  -0.04%  '
$ws.Range("D23").Value = '10.40'
$ws.Range("E23").Value = '  -0.88%  '
$ws.Range("D24").Value = '6.423'
$ws.Range("E24").Value = '  -0.51%  '
$ws.Range("E25").Value = '  +0.05%  '
$ws.Range("D26").Value = '1.787'
$ws.Range("E26").Value = '  -2.94%  '
$ws.Range("E27").Value = '  -0.25%  '
$ws.Range("D28").Value = '2.079'
$ws.Range("E28").Value = '  -1.04%  '
$ws.Range("D29").Value = '113.41'
$ws.Range("E29").Value = '  +0.63%  '
$ws.Range("D30").Value = '4.709'
$ws.Range("E30").Value = '  +0.73%  '
$ws.Range("D31").Value = '4.692'
$ws.Range("E31").Value = '  -0.41%  '
$ws.Range("D32").Value = '0.09153'
$ws.Range("E32").Value = '  +1.20%  '
$ws.Range("D33").Value = '0.05145'
$ws.Range("E33").Value = '  +0.19%  '
$ws.Range("D34").Value = '0.7520'
$ws.Range("E34").Value = '  +3.23%  '
$ws.Range("D35").Value = '2.987'
$ws.Range("E35").Value = '  -1.68%  '
$ws.Range("E36").Value = '  +0.67%  '
$ws.Range("D37").Value = '3.227'
$ws.Range("E37").Value = '  +6.29%  '
$ws.Range("D38").Value = '2.566'
$ws.Range("E38").Value = '  +3.98%  '
$ws.Range("D39").Value = '0.5658'
$ws.Range("E39").Value = '  +6.96%  '
$ws.Range("E40").Value = '  -1.78%  '
$ws.Range("E41").Value = '  -0.19%  '
$ws.Range("D42").Value = '6.615'
$ws.Range("E42").Value = '  +1.22%  '
$ws.Range("D43").Value = '115.66'
$ws.Range("E43").Value = '  -0.09%  '
$ws.Range("D44").Value = '8.532'
$ws.Range("E44").Value = '  +3.19%  '
$ws.Range("D45").Value = '0.1476'
$ws.Range("E45").Value = '  +0.45%  '
$ws.Range("D46").Value = '0.4737'
$ws.Range("E46").Value = '  +2.60%  '
$ws.Range("B47").Value = 'EnergySwap'
$ws.Range("C47").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D47").Value = '10.17'
$ws.Range("E47").Value = '  +1.81%  '
$ws.Range("B48").Value = 'PaxDollar'
$ws.Range("C48").Value = 'https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp'
$ws.Range("D48").Value = '1.000'
$ws.Range("E48").Value = '  +0.19%  '
$ws.Range("D49").Value = '1.566'
$ws.Range("E49").Value = '  +0.13%  '
$ws.Range("D50").Value = '36.91'
$ws.Range("E50").Value = '  +1.07%  '
$ws.Range("D51").Value = '63.13'
$ws.Range("E51").Value = '  -1.15%  '

# Restore default style on the text-forced cells
$ws.Range("D5").Style = "Normal"
$ws.Range("D7").Style = "Normal"
$ws.Range("D8").Style = "Normal"
$ws.Range("D9").Style = "Normal"
$ws.Range("D10").Style = "Normal"
$ws.Range("D11").Style = "Normal"
$ws.Range("D12").Style = "Normal"
$ws.Range("D14").Style = "Normal"
$ws.Range("D16").Style = "Normal"
$ws.Range("D17").Style = "Normal"
$ws.Range("D18").Style = "Normal"
$ws.Range("D19").Style = "Normal"
$ws.Range("D23").Style = "Normal"
$ws.Range("D24").Style = "Normal"
$ws.Range("D26").Style = "Normal"
$ws.Range("D28").Style = "Normal"
$ws.Range("D29").Style = "Normal"
$ws.Range("D30").Style = "Normal"
$ws.Range("D31").Style = "Normal"
$ws.Range("D32").Style = "Normal"
$ws.Range("D33").Style = "Normal"
$ws.Range("D34").Style = "Normal"
$ws.Range("D35").Style = "Normal"
$ws.Range("D37").Style = "Normal"
$ws.Range("D38").Style = "Normal"
$ws.Range("D39").Style = "Normal"
$ws.Range("D42").Style = "Normal"
$ws.Range("D43").Style = "Normal"
$ws.Range("D44").Style = "Normal"
$ws.Range("D45").Style = "Normal"
$ws.Range("D46").Style = "Normal"
$ws.Range("D47").Style = "Normal"
$ws.Range("D48").Style = "Normal"
$ws.Range("D49").Style = "Normal"
$ws.Range("D50").Style = "Normal"
$ws.Range("D51").Style = "Normal"

